# "update 3rd round teaching" — 3rd round edit pass on the deck.
#
# Slide 10 (SlideID 292) title "Scope of identifiers" is renamed to the
# more accurate "Scope of variable names" (shape Id 2, the "Title 1"
# placeholder).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$titleShape = $s.Shapes.Item(2)
$titleShape.TextFrame.TextRange.Text = "Scope of variable names"
